{"js": "// Replace the multiplication-problem text in each table cell with the\n// updated operands, matching the author's regenerated worksheet output.\n// Each old expression is unique in the document, so an exact, case-sensitive\n// search-and-replace on the run text is unambiguous.\nconst replacements = [\n  [\"786\u00d73=\", \"486\u00d77=\"],\n  [\"610\u00d73=\", \"520\u00d75=\"],\n  [\"786\u00d74=\", \"667\u00d76=\"],\n  [\"830\u00d79=\", \"953\u00d76=\"],\n  [\"224\u00d75=\", \"225\u00d75=\"],\n  [\"465\u00d77=\", \"965\u00d79=\"],\n  [\"856\u00d74=\", \"920\u00d76=\"],\n  [\"791\u00d72=\", \"684\u00d76=\"],\n  [\"256\u00d75=\", \"310\u00d79=\"],\n  [\"300\u00d74=\", \"972\u00d77=\"],\n  [\"722\u00d74=\", \"852\u00d78=\"],\n  [\"607\u00d78=\", \"386\u00d74=\"],\n  [\"910\u00d77=\", \"245\u00d77=\"],\n  [\"348\u00d75=\", \"894\u00d76=\"],\n  [\"729\u00d75=\", \"139\u00d74=\"],\n  [\"590\u00d78=\", \"623\u00d79=\"],\n  [\"253\u00d75=\", \"121\u00d75=\"],\n  [\"409\u00d78=\", \"177\u00d78=\"],\n  [\"227\u00d72=\", \"690\u00d79=\"],\n  [\"735\u00d79=\", \"373\u00d78=\"],\n  [\"945\u00d72=\", \"975\u00d72=\"],\n  [\"632\u00d78=\", \"810\u00d76=\"],\n  [\"603\u00d73=\", \"345\u00d76=\"],\n  [\"552\u00d72=\", \"864\u00d78=\"],\n  [\"737\u00d79=\", \"557\u00d75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the multiplication-problem text in each table cell with the\n# updated operands, matching the author's regenerated worksheet output.\n# Each old expression is unique in the document, so Find/Replace-All on the\n# exact (case-sensitive) text is unambiguous and touches only that one run.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"786\u00d73=\", \"486\u00d77=\"),\n    @(\"610\u00d73=\", \"520\u00d75=\"),\n    @(\"786\u00d74=\", \"667\u00d76=\"),\n    @(\"830\u00d79=\", \"953\u00d76=\"),\n    @(\"224\u00d75=\", \"225\u00d75=\"),\n    @(\"465\u00d77=\", \"965\u00d79=\"),\n    @(\"856\u00d74=\", \"920\u00d76=\"),\n    @(\"791\u00d72=\", \"684\u00d76=\"),\n    @(\"256\u00d75=\", \"310\u00d79=\"),\n    @(\"300\u00d74=\", \"972\u00d77=\"),\n    @(\"722\u00d74=\", \"852\u00d78=\"),\n    @(\"607\u00d78=\", \"386\u00d74=\"),\n    @(\"910\u00d77=\", \"245\u00d77=\"),\n    @(\"348\u00d75=\", \"894\u00d76=\"),\n    @(\"729\u00d75=\", \"139\u00d74=\"),\n    @(\"590\u00d78=\", \"623\u00d79=\"),\n    @(\"253\u00d75=\", \"121\u00d75=\"),\n    @(\"409\u00d78=\", \"177\u00d78=\"),\n    @(\"227\u00d72=\", \"690\u00d79=\"),\n    @(\"735\u00d79=\", \"373\u00d78=\"),\n    @(\"945\u00d72=\", \"975\u00d72=\"),\n    @(\"632\u00d78=\", \"810\u00d76=\"),\n    @(\"603\u00d73=\", \"345\u00d76=\"),\n    @(\"552\u00d72=\", \"864\u00d78=\"),\n    @(\"737\u00d79=\", \"557\u00d75=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute(\n        [ref]$oldText,  # FindText\n        [ref]$true,     # MatchCase\n        [ref]$false,    # MatchWholeWord\n        [ref]$false,    # MatchWildcards\n        [ref]$false,    # MatchSoundsLike\n        [ref]$false,    # MatchAllWordForms\n        [ref]$true,     # Forward\n        [ref]1,         # Wrap = wdFindContinue\n        [ref]$false,    # Format\n        [ref]$newText,  # ReplaceWith\n        [ref]2          # Replace = wdReplaceAll\n    )\n}\n"}
